# Fruta / hortaliza, semanal
# Update Fecha (D), Volumen (M), Precio minimo (N), Precio maximo (O),
# Precio promedio ponderado (P) and Precio $/Kg (S) for each data row
# (2-20) of the active sheet with the new weekly values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
    2  = @{ D = 44377; M = 25; N = 20000; O = 20000; P = 20000; S = 1000 }
    3  = @{ D = 44406; M = 20; N = 20000; O = 20000; P = 20000; S = 1000 }
    4  = @{ D = 44291; M = 70; N = 25000; O = 25000; P = 25000; S = 1250 }
    5  = @{ D = 44301; M = 38; N = 22000; O = 22000; P = 22000; S = 1100 }
    6  = @{ D = 44298; M = 65; N = 22000; O = 22000; P = 22000; S = 1100 }
    7  = @{ D = 44385; M = 36; N = 20000; O = 20000; P = 20000; S = 1000 }
    8  = @{ D = 44445; M = 45; N = 20000; O = 20000; P = 20000; S = 1000 }
    9  = @{ D = 44305; M = 20; N = 22000; O = 22000; P = 22000; S = 1100 }
    10 = @{ D = 44376; M = 38; N = 20000; O = 20000; P = 20000; S = 1000 }
    11 = @{ D = 44413; M = 45; N = 20000; O = 20000; P = 20000; S = 1000 }
    12 = @{ D = 44389; M = 20; N = 20000; O = 20000; P = 20000; S = 1000 }
    13 = @{ D = 44400; M = 45; N = 20000; O = 20000; P = 20000; S = 1000 }
    14 = @{ D = 44292; M = 30; N = 25000; O = 25000; P = 25000; S = 1250 }
    15 = @{ D = 44300; M = 45; N = 22000; O = 22000; P = 22000; S = 1100 }
    16 = @{ D = 44448; M = 30; N = 22000; O = 22000; P = 22000; S = 1100 }
    17 = @{ D = 44307; M = 30; N = 22000; O = 22000; P = 22000; S = 1100 }
    18 = @{ D = 44403; M = 50; N = 20000; O = 20000; P = 20000; S = 1000 }
    19 = @{ D = 44294; M = 25; N = 25000; O = 25000; P = 25000; S = 1250 }
    20 = @{ D = 44382; M = 24; N = 20000; O = 20000; P = 20000; S = 1000 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $vals.D   # D: Fecha
    $ws.Cells.Item($r, 13).Value = $vals.M   # M: Volumen
    $ws.Cells.Item($r, 14).Value = $vals.N   # N: Precio minimo
    $ws.Cells.Item($r, 15).Value = $vals.O   # O: Precio maximo
    $ws.Cells.Item($r, 16).Value = $vals.P   # P: Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $vals.S   # S: Precio $/Kg
}
